$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "A11,A12"
$ws.Range("E2").Value = "65,65"

$ws.Range("C3").Value = "A1"
$ws.Range("C4").Value = "A2"
$ws.Range("C5").Value = "A3"
$ws.Range("C6").Value = "A4"
$ws.Range("C7").Value = "A5"
$ws.Range("C8").Value = "A6"
$ws.Range("C9").Value = "A7"
$ws.Range("C10").Value = "A8"
$ws.Range("C11").Value = "A9"
$ws.Range("C12").Value = "A10"

$ws.Range("C13").Select()
